$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2.770546300948285

$ws.Range("B3").Value = 0.3464964993005633
$ws.Range("C3").Value = 0.004309184025731883
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.003598172418533

$ws.Range("B4").Value = 0.006876353814593728
$ws.Range("C4").Value = 0.05231270169004087
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.7119815445968727
